# Generate Report for Archive
#
# 1. Change the status text "Ready for handoff" -> "In Translation" everywhere
#    it appears (Overview sheet's zh-cn/de-de columns, and the Status column
#    on each per-language sheet).
# 2. Narrow the now-shorter status columns (Overview!E:F, zh-cn!C, de-de!C)
#    to match the new text's width.

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"

# NOTE: this host's Range.Value *getter* is unreliable (returns a property
# description instead of the cell's data), so read/write via Value2 instead,
# and always cast to [string] before comparing (Value2 for a boolean-looking
# string like "True" comes back as an actual PowerShell boolean, which would
# otherwise make "-eq $oldText" compare as boolean true/false).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $current = [string]$cell.Value2
        if ($current -eq $oldText) {
            $cell.Value2 = $newText
        }
    }
}

# Narrow the status columns that held the old, longer text.
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E1").EntireColumn.ColumnWidth = 12.5
$ov.Range("F1").EntireColumn.ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
